# DevTesting_IC.wet.xlsx edits:
#  - Edited execution constraints / initialization ruleset to fire for the
#    development testing tool: the active tab moves from "Reservoirs" to
#    "CoordinatedOps" (book-level activeTab + per-sheet tabSelected/selection).
#  - Edited IC to include release tier inputs: ReleaseTier_Input (column C)
#    on CoordinatedOps rows 3-5 is changed from the text "NaN" to the
#    numeric value 0.

$wb = $excel.ActiveWorkbook

# CoordinatedOps holds the release-tier IC inputs and becomes the active sheet.
$ws = $wb.Worksheets.Item("CoordinatedOps")

# Give CoordinatedOps release tier inputs: replace the "NaN" text markers
# with a numeric 0 in column C (MTOMRunType.ReleaseTier_Input) for rows 3-5.
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("C5").Value = 0

# Make CoordinatedOps the active/selected sheet (was Reservoirs) and leave
# the selection on F15, matching the saved view state.
$ws.Activate()
[void]$ws.Range("F15").Select()
